$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ITI")
$ws.Activate()
$src = $ws.Range("I2")
$dst = $ws.Range("I2:I22")
$src.Copy()
$dst.PasteSpecial(-4122)
$dst.Value = "Resolvido"
$dst.Interior.Color = 65535
$ws.Range("I21:I22").Select()
